$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell 'D2' '25.858.18'
Set-TextCell 'D3' '1.738.65'
Set-TextCell 'D4' '0.9999'
Set-TextCell 'E4' '  +0.02%  '
Set-TextCell 'D5' '237.72'
Set-TextCell 'E5' '  +3.16%  '
Set-TextCell 'D6' '0.9998'
Set-TextCell 'E6' '  -0.01%  '
Set-TextCell 'D7' '0.5142'
Set-TextCell 'E7' '  -1.12%  '
Set-TextCell 'D8' '0.2730'
Set-TextCell 'E8' '  -0.66%  '
Set-TextCell 'B9' 'OKB'
Set-TextCell 'C9' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D9' '40.03'
Set-TextCell 'E9' '  +1.68%  '
Set-TextCell 'B10' 'Dogecoin'
Set-TextCell 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 'D10' '0.06123'
Set-TextCell 'E10' '  -0.07%  '
Set-TextCell 'B11' 'WrappedEther'
Set-TextCell 'C11' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D11' '1.738.72'
Set-TextCell 'E11' '  +0.06%  '
Set-TextCell 'B12' 'TRON'
Set-TextCell 'C12' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 'D12' '0.07174'
Set-TextCell 'E12' '  +1.99%  '
Set-TextCell 'D13' '0.6428'
Set-TextCell 'E13' '  +1.37%  '
Set-TextCell 'B14' 'Solana'
Set-TextCell 'C14' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 'D14' '14.89'
Set-TextCell 'E14' '  -0.45%  '
Set-TextCell 'B15' 'Polkadot'
Set-TextCell 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D15' '4.589'
Set-TextCell 'E15' '  +1.62%  '
Set-TextCell 'B16' 'Litecoin'
Set-TextCell 'C16' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D16' '77.21'
Set-TextCell 'E16' '  +0.85%  '
Set-TextCell 'B17' 'Dai'
Set-TextCell 'C17' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D17' '0.9994'
Set-TextCell 'E17' '  -0.04%  '
Set-TextCell 'B18' 'BinanceUSD'
Set-TextCell 'C18' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D18' '0.9999'
Set-TextCell 'E18' '  +0.04%  '
Set-TextCell 'B19' 'WrappedBTC'
Set-TextCell 'C19' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 'D19' '25.874.27'
Set-TextCell 'E19' '  +0.26%  '
Set-TextCell 'B20' 'Avalanche'
Set-TextCell 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 'D20' '11.71'
Set-TextCell 'E20' '  +2.30%  '
Set-TextCell 'B21' 'ShibaInu'
Set-TextCell 'C21' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D21' '0.000006753'
Set-TextCell 'E21' '  +1.88%  '
Set-TextCell 'B22' 'WrappedliquidstakedEther2.0'
Set-TextCell 'C22' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 'D22' '1.962.30'
Set-TextCell 'E22' '  +0.24%  '
Set-TextCell 'B23' 'Uniswap'
Set-TextCell 'C23' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D23' '4.254'
Set-TextCell 'E23' '  +1.88%  '
Set-TextCell 'B24' 'Cosmos'
Set-TextCell 'C24' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D24' '8.662'
Set-TextCell 'E24' '  -1.07%  '
Set-TextCell 'B25' 'Chainlink'
Set-TextCell 'C25' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 'D25' '5.229'
Set-TextCell 'E25' '  +1.93%  '
Set-TextCell 'B26' 'Monero'
Set-TextCell 'C26' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D26' '138.58'
Set-TextCell 'E26' '  -0.48%  '
Set-TextCell 'B27' 'Toncoin'
Set-TextCell 'C27' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D27' '1.510'
Set-TextCell 'E27' '  +0.51%  '
Set-TextCell 'B28' 'EthereumClassic'
Set-TextCell 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D28' '15.14'
Set-TextCell 'E28' '  +1.06%  '
Set-TextCell 'B29' 'LidoDAOToken'
Set-TextCell 'C29' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 'D29' '1.755'
Set-TextCell 'E29' '  -1.08%  '
Set-TextCell 'B30' 'BitcoinCash'
Set-TextCell 'C30' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 'D30' '105.72'
Set-TextCell 'E30' '  +3.81%  '
Set-TextCell 'B31' 'InternetComputer(DFINITY)'
Set-TextCell 'C31' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D31' '3.995'
Set-TextCell 'E31' '  +8.11%  '
Set-TextCell 'B32' 'Stellar'
Set-TextCell 'C32' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D32' '0.08304'
Set-TextCell 'E32' '  +0.60%  '
Set-TextCell 'B33' 'Filecoin'
Set-TextCell 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D33' '3.636'
Set-TextCell 'E33' '  +4.04%  '
Set-TextCell 'B34' 'Hedera'
Set-TextCell 'C34' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D34' '0.04583'
Set-TextCell 'E34' '  +2.71%  '
Set-TextCell 'B35' 'HuobiToken'
Set-TextCell 'C35' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D35' '2.657'
Set-TextCell 'E35' '  +1.93%  '
Set-TextCell 'B36' 'ARBITRUM'
Set-TextCell 'C36' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D36' '0.9862'
Set-TextCell 'E36' '  +1.54%  '
Set-TextCell 'B37' 'ImmutableX'
Set-TextCell 'C37' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D37' '0.6171'
Set-TextCell 'E37' '  +0.54%  '
Set-TextCell 'B38' 'MXToken'
Set-TextCell 'C38' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D38' '2.685'
Set-TextCell 'E38' '  +0.32%  '
Set-TextCell 'B39' 'VeChain'
Set-TextCell 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D39' '0.01611'
Set-TextCell 'E39' '  +2.75%  '
Set-TextCell 'B40' 'RenderToken'
Set-TextCell 'C40' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D40' '1.926'
Set-TextCell 'E40' '  +1.12%  '
Set-TextCell 'B41' 'PaxDollar'
Set-TextCell 'C41' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D41' '0.9994'
Set-TextCell 'E41' '  -0.02%  '
Set-TextCell 'B42' 'Quant'
Set-TextCell 'C42' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D42' '97.49'
Set-TextCell 'E42' '  -2.29%  '
Set-TextCell 'B43' 'TheSandbox'
Set-TextCell 'C43' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 'D43' '0.3828'
Set-TextCell 'E43' '  +0.11%  '
Set-TextCell 'B44' 'TrustWalletToken'
Set-TextCell 'C44' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D44' '0.7359'
Set-TextCell 'E44' '  +1.60%  '
Set-TextCell 'B45' 'FraxShare'
Set-TextCell 'C45' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D45' '4.947'
Set-TextCell 'E45' '  -0.87%  '
Set-TextCell 'B46' 'Algorand'
Set-TextCell 'C46' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D46' '0.1124'
Set-TextCell 'E46' '  -0.15%  '
Set-TextCell 'B47' 'Cronos'
Set-TextCell 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D47' '0.05256'
Set-TextCell 'E47' '  -1.66%  '
Set-TextCell 'B48' 'Aptos'
Set-TextCell 'C48' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D48' '6.147'
Set-TextCell 'E48' '  -0.04%  '
Set-TextCell 'B49' 'Aave'
Set-TextCell 'C49' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D49' '54.79'
Set-TextCell 'E49' '  +3.34%  '
Set-TextCell 'B50' 'Elrond'
Set-TextCell 'C50' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 'D50' '30.43'
Set-TextCell 'E50' '  +1.69%  '
Set-TextCell 'B51' 'EnergySwap'
Set-TextCell 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D51' '7.575'
Set-TextCell 'E51' '  +0.04%  '
